$d = $word.ActiveDocument

# Locate the last two paragraphs of the document body:
#   Para A (second to last) -- holds the (hidden) "_GoBack" bookmark, styled
#            with Consolas/bold/green/sz18.
#   Para B (last)           -- a trailing empty paragraph, Consolas/sz18.
$n = $d.Paragraphs.Count
$paraA = $d.Paragraphs.Item($n - 1)
$paraB = $d.Paragraphs.Item($n)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Step 1: replace the trailing empty paragraph (Para B) with the whole
# new block of content (everything that follows Para A in the diff). Doing
# this before touching Para A keeps Para A's Range valid/addressable.
$targetB = $d.Range($paraB.Range.Start, $paraB.Range.End)

$xmlB = @"
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t>Resultado Final: Implementação e testes realizados com sucesso em aplicações de cenários positivos e negativos na simulação</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">Gabriel Siqueira </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
    <w:t>Petillo</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
    <w:t>RM 81238 – 3SIA 2020</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
    </w:rPr>
    <w:t>FIAP</w:t>
  </w:r>
</w:p>
"@

$targetB.InsertXML($xmlB)

# --- Step 2: strip the bookmark out of Para A, turning it into a plain
# empty paragraph (same pPr/rPr it always had). Para A's Range is still
# valid since step 1 only touched content after it.
$targetA = $d.Range($paraA.Range.Start, $paraA.Range.End)

$xmlA = @"
<w:p $w>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
      <w:b/>
      <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$targetA.InsertXML($xmlA)

# --- Step 3: re-home the (single, document-unique) "_GoBack" bookmark onto
# the new "signature block" -- from the start of the centered paragraph
# that used to carry bookmarkStart through to right after "FIAP" (where
# bookmarkEnd used to sit), mirroring the diff exactly.
$total = $d.Paragraphs.Count
$fiapPara = $d.Paragraphs.Item($total)
$startPara = $d.Paragraphs.Item($total - 3)

$bmRange = $d.Range($startPara.Range.Start, $fiapPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
